$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3").Value = "Test Baltic Spirit"
$ws.Range("A3").Value = "TBSJ"
$ws.Range("B3").Select() | Out-Null
